# V 2.0.2 se arreglo la fecha y hora de reimpresion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient identification block ---
$ws.Range("A6").Value = "XX  XX MASCULINO     "
$ws.Range("G6").Value = "/201762652"

# Fecha de nacimiento (text that looks like a date -> force text format
# so it is not auto-converted into a date serial number)
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "1991-10-24"

$ws.Range("D9").Value = "26 AÑOS APROX"
$ws.Range("E9").Value = "CAPITAL"

# --- Estado civil / ocupacion / nacionalidad / documento row ---
$ws.Range("A11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("E11").Value = "GUATEMALTECA"
$ws.Range("G11").Value = ""

# --- Emergencia contact row ---
$ws.Range("A13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = "IGNORADA"
$ws.Range("G13").Value = ""

# --- Fecha / hora de asistencia ---
$ws.Range("D14").Value = "Hora: 15:50:19"

# Fecha (dd/mm/yyyy) - stays as text automatically since this engine does not
# parse dd/mm/yyyy as a date by default (unlike the yyyy-mm-dd case above)
$ws.Range("A15").Value = "24/10/2017"

# --- Tipo de consulta ---
$ws.Range("D16").Value = ""

Write-Host "edits applied"
